$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Binomial")

# Fill the "Percentage of Outcomes" column with the binomial distribution
# formula (this reproduces the shared formula B5:B25 from the diff).
$ws.Range("B5:B25").Formula = '=BINOM.DIST(A5,$B$2,$B$1,FALSE)'

# Match the author's final selection (B1 -> B5).
$ws.Range("B5").Select()
